$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes: insert new rows, copying formatting from an existing
# row that uses the same "s=4" style (row 6 / page1-district row) ---

# 1) Insert a new row at position 6 (becomes page1-date_Monday / date_week)
$ws.Rows("6:6").Insert()
$ws.Range("A7:B7").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)

# 2) Insert 3 new rows at position 12 (pox_present / probe_present / pox_func)
$ws.Rows("12:14").Insert()
$ws.Range("A15:B15").Copy()
$ws.Range("A12:B14").PasteSpecial(-4122)

# 3) Insert 5 new rows at position 19 (monday_seen .. friday_seen)
$ws.Rows("19:23").Insert()
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Value changes, written in the exact order needed so the shared-string
# table ends up with the same append order as the target workbook ---

$ws.Range("A6").Value = "page1-date_Monday"
$ws.Range("B6").Value = "date_week"

$ws.Range("A19").Value = "page1-monday_seen"
$ws.Range("A20").Value = "page1-tuesday_seen"
$ws.Range("A21").Value = "page1-wednesday_seen"
$ws.Range("A22").Value = "page1-thursday_seen"
$ws.Range("A23").Value = "page1-friday_seen"

$ws.Range("B18").Value = "u5_attendance_last_week"
$ws.Range("B19").Value = "u5_attendance_last_mon"
$ws.Range("B20").Value = "u5_attendance_last_tue"
$ws.Range("B21").Value = "u5_attendance_last_wed"
$ws.Range("B22").Value = "u5_attendance_last_thu"
$ws.Range("B23").Value = "u5_attendance_last_fri"

$ws.Range("A12").Value = "page1-pox_present"
$ws.Range("A13").Value = "page1-probe_present"
$ws.Range("A14").Value = "page1-pox_func"

$ws.Range("B13").Value = "probe_avail"
$ws.Range("B12").Value = "pox_avail"
$ws.Range("B14").Value = "pox_functional"

# --- Column widths ---
# (ColumnWidth is stored internally rounded to 1/7-character steps plus a
# fixed 5/7 padding offset, so we back-solve for the input that reproduces
# the target stored width as closely as the engine's precision allows)
$ws.Columns("A").ColumnWidth = 20.285714285714285
$ws.Columns("B").ColumnWidth = 21.660714285714285

# --- Selection / active cell ---
$ws.Range("D15").Select()

Write-Host "edit applied"
